$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment keeps inlineStr/shared-string text type
$ws.Range('D2').Value = '62.896.26'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '2.461.74'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '2.461.14'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '2.907.97'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').Value = '62.787.31'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '2.462.28'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('E19').Value = '  +2.67%  '
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  +11.08%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('E25').Value = '  +19.86%  '
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +2.10%  '
$ws.Range('D28').Value = '0.0₃0980'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  -13.19%  '
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('E34').Value = '  -3.51%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').Value = '0.0₆0321'
$ws.Range('E44').Value = '  -74.91%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  +6.72%  '
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  +0.09%  '

# Numeric-looking text values: must force text format so Excel keeps them as strings, not numbers,
# then restore default style so no stray style index is left on the cell.
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '575.71'
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '146.08'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.356'
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '29.08'
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0000178'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.03'
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '327.35'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '65.89'
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '651.54'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '8.02'
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.86'
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.76'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.370'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.40'
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '18.74'
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '151.14'
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.76'
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '154.01'
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.59'
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '20.47'
$c.Style = 'Normal'
